# Fix the typo "Bath" -> "Path" in the cycle SmartArt diagram on slide 14
# ("Administracion del ambiente"). The diagram is hosted in a graphicFrame
# placeholder; one of its cycle nodes currently reads "Bath" and should
# read "Path" instead.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasSmartArt) {
            $nodes = $shape.SmartArt.AllNodes
            for ($ni = 1; $ni -le $nodes.Count; $ni++) {
                $node = $nodes.Item($ni)
                $tr = $node.TextFrame2.TextRange
                if ($tr.Text -eq "Bath") {
                    $tr.Text = "Path"
                }
            }
        }
    }
}
